$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BF2:BF31").NumberFormat = "@"

$ws.Range("AT2").Value = 24
$ws.Range("BF2").Value = "2013-04-11"
$ws.Range("AD3").Value = 9
$ws.Range("AO3").Value = 18
$ws.Range("AP3").Value = 19
$ws.Range("AW3").Value = 12
$ws.Range("BA3").Value = 17
$ws.Range("BF3").Value = "2013-04-11"
$ws.Range("AD4").Value = 9
$ws.Range("AK4").Value = 16
$ws.Range("BF4").Value = "2013-04-11"
$ws.Range("AD5").Value = 9
$ws.Range("AJ5").Value = 18
$ws.Range("BB5").Value = 28
$ws.Range("BF5").Value = "2013-04-11"
$ws.Range("AA6").Value = 19.6
$ws.Range("AB6").Value = 92.8
$ws.Range("AC6").Value = 0.2
$ws.Range("AD6").Value = 29
$ws.Range("AE6").Value = 13
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 27
$ws.Range("AK6").Value = 26
$ws.Range("AN6").Value = 24
$ws.Range("AO6").Value = 20
$ws.Range("AP6").Value = 22
$ws.Range("AT6").Value = 9
$ws.Range("AX6").Value = 15
$ws.Range("BA6").Value = 18
$ws.Range("BB6").Value = 30
$ws.Range("BF6").Value = "2013-04-11"
$ws.Range("D6").Value = 77
$ws.Range("E6").Value = 42
$ws.Range("G6").Value = 0.545
$ws.Range("H6").Value = 48.3
$ws.Range("I6").Value = 35.6
$ws.Range("J6").Value = 81.8
$ws.Range("L6").Value = 5.2
$ws.Range("M6").Value = 15.1
$ws.Range("N6").Value = 0.345
$ws.Range("O6").Value = 16.3
$ws.Range("P6").Value = 21.1
$ws.Range("Q6").Value = 0.773
$ws.Range("R6").Value = 12.7
$ws.Range("S6").Value = 30.5
$ws.Range("T6").Value = 43.2
$ws.Range("V6").Value = 14.3
$ws.Range("AD7").Value = 9
$ws.Range("BF7").Value = "2013-04-11"
$ws.Range("AD8").Value = 9
$ws.Range("AL8").Value = 12
$ws.Range("AO8").Value = 19
$ws.Range("BF8").Value = "2013-04-11"
$ws.Range("AD9").Value = 9
$ws.Range("AT9").Value = 3
$ws.Range("BF9").Value = "2013-04-11"
$ws.Range("AH10").Value = 20
$ws.Range("BF10").Value = "2013-04-11"
$ws.Range("AB11").Value = 100.9
$ws.Range("AC11").Value = 0.9
$ws.Range("AD11").Value = 9
$ws.Range("AH11").Value = 17
$ws.Range("AM11").Value = 12
$ws.Range("AP11").Value = 17
$ws.Range("AQ11").Value = 3
$ws.Range("AT11").Value = 2
$ws.Range("AV11").Value = 26
$ws.Range("AY11").Value = 15
$ws.Range("BF11").Value = "2013-04-11"
$ws.Range("D11").Value = 78
$ws.Range("F11").Value = 33
$ws.Range("G11").Value = 0.577
$ws.Range("I11").Value = 38
$ws.Range("M11").Value = 19.9
$ws.Range("P11").Value = 21.4
$ws.Range("Q11").Value = 0.792
$ws.Range("R11").Value = 10.8
$ws.Range("U11").Value = 22.5
$ws.Range("W11").Value = 6.8
$ws.Range("Z11").Value = 21.4
$ws.Range("AD12").Value = 9
$ws.Range("AF12").Value = 11
$ws.Range("AT12").Value = 7
$ws.Range("AZ12").Value = 16
$ws.Range("BF12").Value = "2013-04-11"
$ws.Range("AD13").Value = 9
$ws.Range("BF13").Value = "2013-04-11"
$ws.Range("AD14").Value = 9
$ws.Range("BF14").Value = "2013-04-11"
$ws.Range("AE15").Value = 13
$ws.Range("AH15").Value = 29
$ws.Range("AX15").Value = 13
$ws.Range("BF15").Value = "2013-04-11"
$ws.Range("AD16").Value = 9
$ws.Range("AJ16").Value = 17
$ws.Range("AN16").Value = 23
$ws.Range("AO16").Value = 17
$ws.Range("AP16").Value = 19
$ws.Range("AZ16").Value = 17
$ws.Range("BF16").Value = "2013-04-11"
$ws.Range("AD17").Value = 9
$ws.Range("AX17").Value = 14
$ws.Range("BF17").Value = "2013-04-11"
$ws.Range("AD18").Value = 9
$ws.Range("AK18").Value = 25
$ws.Range("AM18").Value = 14
$ws.Range("BF18").Value = "2013-04-11"
$ws.Range("AD19").Value = 9
$ws.Range("BF19").Value = "2013-04-11"
$ws.Range("BF20").Value = "2013-04-11"
$ws.Range("AB21").Value = 100
$ws.Range("AC21").Value = 4.4
$ws.Range("AD21").Value = 29
$ws.Range("AF21").Value = 6
$ws.Range("AH21").Value = 30
$ws.Range("AJ21").Value = 19
$ws.Range("AK21").Value = 15
$ws.Range("AP21").Value = 16
$ws.Range("AR21").Value = 19
$ws.Range("AT21").Value = 25
$ws.Range("AW21").Value = 10
$ws.Range("AY21").Value = 2
$ws.Range("AZ21").Value = 15
$ws.Range("BF21").Value = "2013-04-11"
$ws.Range("D21").Value = 77
$ws.Range("F21").Value = 26
$ws.Range("G21").Value = 0.662
$ws.Range("I21").Value = 36.4
$ws.Range("J21").Value = 81.2
$ws.Range("K21").Value = 0.449
$ws.Range("M21").Value = 28.8
$ws.Range("N21").Value = 0.377
$ws.Range("S21").Value = 29.6
$ws.Range("T21").Value = 40.6
$ws.Range("U21").Value = 19.3
$ws.Range("Z21").Value = 20.1
$ws.Range("AB22").Value = 105.8
$ws.Range("AC22").Value = 9.199999999999999
$ws.Range("AD22").Value = 9
$ws.Range("AH22").Value = 12
$ws.Range("AL22").Value = 13
$ws.Range("AU22").Value = 21
$ws.Range("AW22").Value = 11
$ws.Range("AY22").Value = 3
$ws.Range("BF22").Value = "2013-04-11"
$ws.Range("D22").Value = 78
$ws.Range("E22").Value = 57
$ws.Range("G22").Value = 0.731
$ws.Range("I22").Value = 38
$ws.Range("J22").Value = 79
$ws.Range("L22").Value = 7.3
$ws.Range("M22").Value = 19.2
$ws.Range("N22").Value = 0.381
$ws.Range("P22").Value = 27.2
$ws.Range("S22").Value = 33.2
$ws.Range("U22").Value = 21.4
$ws.Range("W22").Value = 8.199999999999999
$ws.Range("BF23").Value = "2013-04-11"
$ws.Range("AD24").Value = 9
$ws.Range("BB24").Value = 29
$ws.Range("BF24").Value = "2013-04-11"
$ws.Range("AH25").Value = 20
$ws.Range("BF25").Value = "2013-04-11"
$ws.Range("AD26").Value = 9
$ws.Range("BF26").Value = "2013-04-11"
$ws.Range("AD27").Value = 9
$ws.Range("AK27").Value = 17
$ws.Range("BF27").Value = "2013-04-11"
$ws.Range("AD28").Value = 9
$ws.Range("AE28").Value = 2
$ws.Range("AG28").Value = 2
$ws.Range("AP28").Value = 21
$ws.Range("AQ28").Value = 4
$ws.Range("AY28").Value = 16
$ws.Range("BF28").Value = "2013-04-11"
$ws.Range("AD29").Value = 9
$ws.Range("AU29").Value = 20
$ws.Range("BF29").Value = "2013-04-11"
$ws.Range("BF30").Value = "2013-04-11"
$ws.Range("AT31").Value = 8
$ws.Range("BF31").Value = "2013-04-11"
